$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 45955
$ws.Range("D8").Value = 161.32
$ws.Range("E8").Value = 158.15
$ws.Range("F8").Value = 168.15
$ws.Range("G8").Value = 158.31
$ws.Range("A9").Value = 45955
$ws.Range("D9").Value = 161.32
$ws.Range("E9").Value = 158.15
$ws.Range("F9").Value = 168.15
$ws.Range("G9").Value = 158.31
$ws.Range("A10").Value = 45955
$ws.Range("D10").Value = 163.72
$ws.Range("E10").Value = 160.49
$ws.Range("F10").Value = 170.49
$ws.Range("G10").Value = 160.96
$ws.Range("A11").Value = 45954
$ws.Range("D11").Value = 159.86000000000001
$ws.Range("E11").Value = 157.69999999999999
$ws.Range("F11").Value = 167.7
$ws.Range("G11").Value = 157.86000000000001
$ws.Range("A12").Value = 45954
$ws.Range("D12").Value = 159.86000000000001
$ws.Range("E12").Value = 157.69999999999999
$ws.Range("F12").Value = 167.7
$ws.Range("G12").Value = 157.86000000000001
$ws.Range("A13").Value = 45954
$ws.Range("D13").Value = 162.16
$ws.Range("E13").Value = 159.91
$ws.Range("F13").Value = 169.91
$ws.Range("G13").Value = 160.38
$ws.Range("A17").Value = 45955
$ws.Range("D17").Value = 167.16
$ws.Range("E17").Value = 163.29
$ws.Range("F17").Value = 173.29
$ws.Range("A18").Value = 45954
$ws.Range("D18").Value = 165.59
$ws.Range("E18").Value = 162.78
$ws.Range("F18").Value = 172.78
$ws.Range("A22").Value = 45955
$ws.Range("D22").Value = 162.44999999999999
$ws.Range("E22").Value = 159.5
$ws.Range("F22").Value = 169.1
$ws.Range("G22").Value = 160.66999999999999
$ws.Range("A23").Value = 45955
$ws.Range("D23").Value = 168.49
$ws.Range("E23").Value = 164.21
$ws.Range("F23").Value = 174.21
$ws.Range("A24").Value = 45955
$ws.Range("D24").Value = 168.29
$ws.Range("E24").Value = 164.41
$ws.Range("F24").Value = 174.41
$ws.Range("A25").Value = 45955
$ws.Range("D25").Value = 169.12
$ws.Range("E25").Value = 163.80000000000001
$ws.Range("F25").Value = 173.8
$ws.Range("G25").Value = 163.63
$ws.Range("A26").Value = 45955
$ws.Range("D26").Value = 167.85
$ws.Range("E26").Value = 165.34
$ws.Range("F26").Value = 175.34
$ws.Range("A27").Value = 45954
$ws.Range("D27").Value = 160.88999999999999
$ws.Range("E27").Value = 158.81
$ws.Range("F27").Value = 168.41
$ws.Range("G27").Value = 159.99
$ws.Range("A28").Value = 45954
$ws.Range("D28").Value = 166.93
$ws.Range("E28").Value = 163.63
$ws.Range("F28").Value = 173.63
$ws.Range("A29").Value = 45954
$ws.Range("D29").Value = 166.73
$ws.Range("E29").Value = 163.84
$ws.Range("F29").Value = 173.84
$ws.Range("A30").Value = 45954
$ws.Range("D30").Value = 167.56
$ws.Range("E30").Value = 163.22999999999999
$ws.Range("F30").Value = 173.23
$ws.Range("G30").Value = 163.05000000000001
$ws.Range("A31").Value = 45954
$ws.Range("D31").Value = 166.28
$ws.Range("E31").Value = 164.76
$ws.Range("F31").Value = 174.76
$ws.Range("A35").Value = 45955
$ws.Range("D35").Value = 161.97
$ws.Range("E35").Value = 157.69999999999999
$ws.Range("F35").Value = 166.71
$ws.Range("A36").Value = 45954
$ws.Range("D36").Value = 160.4
$ws.Range("E36").Value = 157.13
$ws.Range("F36").Value = 166.13
$ws.Range("A40").Value = 45955
$ws.Range("D40").Value = 167.61
$ws.Range("E40").Value = 163.01
$ws.Range("F40").Value = 173.01
$ws.Range("A41").Value = 45955
$ws.Range("D41").Value = 167.33
$ws.Range("E41").Value = 163.43
$ws.Range("F41").Value = 173.43
$ws.Range("A42").Value = 45954
$ws.Range("D42").Value = 166.05
$ws.Range("E42").Value = 162.51
$ws.Range("F42").Value = 172.51
$ws.Range("A43").Value = 45954
$ws.Range("D43").Value = 165.77
$ws.Range("E43").Value = 162.93
$ws.Range("F43").Value = 172.93
$ws.Range("A47").Value = 45955
$ws.Range("D47").Value = 160.61000000000001
$ws.Range("E47").Value = 159.08000000000001
$ws.Range("F47").Value = 169.08
$ws.Range("A48").Value = 45955
$ws.Range("D48").Value = 160.59
$ws.Range("E48").Value = 159.25
$ws.Range("F48").Value = 169.25
$ws.Range("A49").Value = 45954
$ws.Range("D49").Value = 159.88
$ws.Range("E49").Value = 158.61000000000001
$ws.Range("F49").Value = 168.61
$ws.Range("A50").Value = 45954
$ws.Range("D50").Value = 159.86000000000001
$ws.Range("E50").Value = 158.78
$ws.Range("F50").Value = 168.78
$ws.Range("A54").Value = 45955
$ws.Range("D54").Value = 177.81
$ws.Range("E54").Value = 173.44
$ws.Range("F54").Value = 183.44
$ws.Range("A55").Value = 45955
$ws.Range("D55").Value = 165.46
$ws.Range("E55").Value = 170.77
$ws.Range("F55").Value = 180.77
$ws.Range("A56").Value = 45955
$ws.Range("D56").Value = 167.74
$ws.Range("A57").Value = 45955
$ws.Range("D57").Value = 167.42
$ws.Range("E57").Value = 165.04
$ws.Range("A58").Value = 45955
$ws.Range("D58").Value = 163.32
$ws.Range("E58").Value = 161.09
$ws.Range("F58").Value = 171.09
$ws.Range("A59").Value = 45955
$ws.Range("D59").Value = 170.15
$ws.Range("E59").Value = 171.65
$ws.Range("A60").Value = 45954
$ws.Range("D60").Value = 176.24
$ws.Range("E60").Value = 172.88
$ws.Range("F60").Value = 182.88
$ws.Range("A61").Value = 45954
$ws.Range("D61").Value = 163.9
$ws.Range("E61").Value = 170.3
$ws.Range("F61").Value = 180.3
$ws.Range("A62").Value = 45954
$ws.Range("D62").Value = 166.18
$ws.Range("A63").Value = 45954
$ws.Range("D63").Value = 165.85
$ws.Range("E63").Value = 164.57
$ws.Range("A64").Value = 45954
$ws.Range("D64").Value = 161.76
$ws.Range("E64").Value = 160.62
$ws.Range("F64").Value = 170.62
$ws.Range("A65").Value = 45954
$ws.Range("D65").Value = 168.58
$ws.Range("E65").Value = 171.08
